# "Generate Report for Handback"
# Updates the handoff/handback timestamps produced by a new handback run
# for the file "0ad0c369-1e2e-4808-9f44-53f2bda72a7c" on the Overview,
# zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet --------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-01 15:08:53"

# --- zh-cn sheet -------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-01 15:08:48"
$wsZhCn.Range("K2").Value = "2016-09-01 15:09:24"

# --- de-de sheet -------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-09-01 15:08:53"
$wsDeDe.Range("K2").Value = "2016-09-01 15:09:32"
